$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Remove the "Confluence - Documentation" bullet paragraph from the
#     "Content Placeholder 7" shape. ---
$content = $s.Shapes.Item("Content Placeholder 7")
$tr = $content.TextFrame.TextRange
$count = $tr.Paragraphs().Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i)
    $txt = $para.Text.TrimEnd("`r", "`n")
    if ($txt -eq "Confluence - Documentation") {
        $para.Delete()
    }
}

# --- Nudge three of the stack-logo pictures down/right a bit. Left/Top are
#     expressed in points (Single precision, like real PowerPoint COM), so
#     the literals below are chosen so they round-trip to the exact target
#     EMU offsets. ---
$pic2 = $s.Shapes.Item("Picture 2")
$pic2.Left = 54.18893051147461
$pic2.Top = 133.59664916992188

$pic3 = $s.Shapes.Item("Picture 3")
$pic3.Left = 75.87508392333984
$pic3.Top = 181.92184448242188

$pic4 = $s.Shapes.Item("Picture 4")
$pic4.Left = 71.15035247802734
$pic4.Top = 255.97398376464844

# --- Drop the last stack-logo picture entirely. ---
$pic5 = $s.Shapes.Item("Picture 5")
$pic5.Delete()
